$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$timestamps = @(
    "2024-11-01 08:09:26",
    "2024-11-01 08:25:22",
    "2024-11-01 08:29:34",
    "2024-11-01 08:31:19",
    "2024-11-01 08:32:26",
    "2024-11-01 08:34:26",
    "2024-11-01 08:42:30",
    "2024-11-01 08:45:54",
    "2024-11-01 08:59:50",
    "2024-11-01 09:24:40",
    "2024-11-01 10:01:44"
)

$startRow = 6
for ($i = 0; $i -lt $timestamps.Length; $i++) {
    $row = $startRow + $i
    $ws.Cells.Item($row, 1).Value = $timestamps[$i]
    $ws.Cells.Item($row, 2).Value = "Success"
}
